$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.487.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.663.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.10%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "646.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.23%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.144"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.440"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000229"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.82%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.282.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.678.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.516.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.27%  "

$ws.Range("E17").Value = "  +0.89%  "

$ws.Range("E18").Value = "  -0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "470.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.648"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.812.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.98%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000125"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.16%  "

$ws.Range("E31").Value = "  -0.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.165"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.662.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.45%  "

$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "179.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.83%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0890"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.98%  "

$ws.Range("E43").Value = "  -2.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.924"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.28%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.85%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.41%  "

$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.09%  "

$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000265"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.91%  "
